# Updates the "cryptos" price list to the latest scraped snapshot.
# Mirrors the GitHub Actions job that refreshes cryptos.xlsx: for each
# changed row, the Price (D) and Volume(1h) (E) columns are refreshed,
# and a couple of rows (19/20, 30/31, 42/43) swapped rank order, which
# also updates their Coin (B) and Link (C) columns.
#
# Many Price values are numeric-looking strings ("1.00", "209.46", ...).
# Assigning those to Range.Value would make Excel auto-convert them to
# real numbers (dropping the trailing zero / formatting), so - exactly
# like typing an apostrophe-prefixed entry in Excel - those are written
# with a leading ' to force text storage while keeping the visible text
# identical. Values that are not valid numbers (e.g. "27.040.29",
# "1.561.59") already round-trip as text with no special handling.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '27.040.29'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '1.561.59'
$ws.Range('E3').Value = '  +0.64%  '
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.21%  '
$ws.Range('D5').Value = '''209.46'
$ws.Range('E5').Value = '  +1.51%  '
$ws.Range('D6').Value = '''0.490'
$ws.Range('E6').Value = '  +0.41%  '
$ws.Range('E7').Value = '  +0.52%  '
$ws.Range('D8').Value = '''21.95'
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('E9').Value = '  -0.17%  '
$ws.Range('D10').Value = '''0.0595'
$ws.Range('E10').Value = '  +0.28%  '
$ws.Range('D11').Value = '''0.0861'
$ws.Range('E11').Value = '  +0.55%  '
$ws.Range('D12').Value = '1.785.32'
$ws.Range('E12').Value = '  +0.71%  '
$ws.Range('D13').Value = '1.561.64'
$ws.Range('E13').Value = '  +0.59%  '
$ws.Range('E14').Value = '  +0.58%  '
$ws.Range('E15').Value = '  +0.00%  '
$ws.Range('D16').Value = '27.039.47'
$ws.Range('E16').Value = '  +0.48%  '
$ws.Range('E17').Value = '  +0.56%  '
$ws.Range('E18').Value = '  -1.17%  '
$ws.Range('B19').Value = 'BitcoinCash'
$ws.Range('C19').Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range('D19').Value = '''215.33'
$ws.Range('E19').Value = '  -0.71%  '
$ws.Range('B20').Value = 'Chainlink'
$ws.Range('C20').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D20').Value = '''7.39'
$ws.Range('E20').Value = '  +1.39%  '
$ws.Range('E21').Value = '  +0.42%  '
$ws.Range('D22').Value = '''4.13'
$ws.Range('E22').Value = '  +1.20%  '
$ws.Range('D23').Value = '''9.15'
$ws.Range('E23').Value = '  -0.29%  '
$ws.Range('D24').Value = '''1.93'
$ws.Range('E24').Value = '  -0.51%  '
$ws.Range('D25').Value = '''153.79'
$ws.Range('E25').Value = '  +0.10%  '
$ws.Range('D26').Value = '''6.59'
$ws.Range('E26').Value = '  -0.82%  '
$ws.Range('E27').Value = '  +0.68%  '
$ws.Range('E28').Value = '  +0.98%  '
$ws.Range('E29').Value = '  +0.46%  '
$ws.Range('B30').Value = 'PancakeSwap'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D30').Value = '''1.13'
$ws.Range('E30').Value = '  +4.27%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '''0.0471'
$ws.Range('E31').Value = '  +0.67%  '
$ws.Range('E32').Value = '  +0.44%  '
$ws.Range('D33').Value = '''3.18'
$ws.Range('E33').Value = '  +2.21%  '
$ws.Range('D34').Value = '1.428.87'
$ws.Range('E34').Value = '  +1.36%  '
$ws.Range('E35').Value = '  +16.25%  '
$ws.Range('E36').Value = '  +0.88%  '
$ws.Range('E37').Value = '  +2.69%  '
$ws.Range('E38').Value = '  +0.79%  '
$ws.Range('D39').Value = '''0.531'
$ws.Range('E39').Value = '  +1.32%  '
$ws.Range('D40').Value = '''5.84'
$ws.Range('E40').Value = '  +3.44%  '
$ws.Range('D41').Value = '''0.807'
$ws.Range('E41').Value = '  +0.05%  '
$ws.Range('B42').Value = 'PaxDollar'
$ws.Range('C42').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D42').Value = '''1.01'
$ws.Range('E42').Value = '  +0.47%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D43').Value = '''2.36'
$ws.Range('E43').Value = '  +2.67%  '
$ws.Range('D44').Value = '''0.999'
$ws.Range('E44').Value = '  +0.31%  '
$ws.Range('D45').Value = '''64.26'
$ws.Range('E45').Value = '  -0.29%  '
$ws.Range('D46').Value = '''1.74'
$ws.Range('E46').Value = '  +0.64%  '
$ws.Range('D47').Value = '1.702.49'
$ws.Range('E47').Value = '  +0.94%  '
$ws.Range('E48').Value = '  -2.03%  '
$ws.Range('E49').Value = '  +1.11%  '
$ws.Range('D50').Value = '''0.0515'
$ws.Range('E50').Value = '  +0.28%  '
$ws.Range('E51').Value = '  +0.25%  '
